# Update ticket price / want-to-go counts on the "展览" and "全部类型" sheets.
# G2: 50 -> 51.4 (最低票价 for row 2)
# F6: 797 -> 799 (想去人数 for row 6)

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G2").Value = 51.4
    $ws.Range("F6").Value = 799
}
